$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The PA10 row is being removed (column A shifts up one row starting at A12)
# and the PA6 deadline is corrected (column B shifts up one row starting at B8),
# then a new Studio Assignment 8 row is appended at the end.

# --- Column A: remove "PA10" (old A12), shifting A13:A19 up into A12:A18 ---
$colA = @("PA1","PA2","PA3","PA4","PA5a","PA5b","PA6","PA7","PA8","PA9","SA1","SA2","SA3","SA4","SA5","SA6","SA7","SA8")

# --- Column B: remove "October 23, 2024" (old B8), shifting B9:B19 up into B8:B18 ---
$colB = @("September 11, 2024","September 18, 2024","September 27, 2024","October 04, 2024","October 11, 2024","October 16, 2024","October 30, 2024","November 20, 2024","December 04, 2024","December 11, 2024","September 09, 2024","September 16, 2024","September 23, 2024","September 30, 2024","October 07, 2024","October 09, 2024","October 21, 2024","October 28, 2024")

for ($i = 0; $i -lt $colA.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $colA[$i]
    $ws.Cells.Item($row, 2).Value = $colB[$i]
}

# The brand-new row (19) needs the same text number format the rest of column B uses
$ws.Cells.Item(19, 2).NumberFormat = "@"

# Update selection to match the end-state cursor position (just past the new last row)
$ws.Range("B20").Select()
